# Update cryptos list cell values per the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "215.00"); force
# text format first so Excel does not silently coerce them to numbers
# and drop significant trailing zeros.
$ws.Range('D2').Value = '26.950.89'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.674.34'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.00'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.26'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '1.910.82'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').Value = '1.658.79'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.526'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.63'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '26.966.02'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.14'
$ws.Range('E18').Value = '  +4.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '234.86'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.44'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.18'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('E24').Value = '  -2.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.75'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.18'
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.03'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').Value = '1.487.57'
$ws.Range('E33').Value = '  -4.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.68'
$ws.Range('E35').Value = '  +3.29%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.583'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.897'
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('E40').Value = '  +8.16%  '
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.30'
$ws.Range('E43').Value = '  +2.57%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '67.46'
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('D45').Value = '1.818.52'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.60'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0508'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.70'
$ws.Range('E51').Value = '  +0.07%  '
